# Botones de retroceso en usuario agregado
# Extends the "usuario" sheet data with four more placeholder rows
# (Id 10-13) so the "back"/pagination buttons have rows to page through.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 11
$startId = 10
$rowsToAdd = 4

for ($i = 0; $i -lt $rowsToAdd; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $startId + $i

    # Columns B:E stay empty text placeholders, matching the existing
    # blank rows (e.g. row 10) already present in the sheet.
    $blankRange = $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 5))
    $blankRange.Value = "'"
    $blankRange.Style = "Normal"
}
